$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 926338.2  # ALC!H17: 1163285 -> 926338.2
$ws.Cells.Item(17, 10).Value = 926338.2  # ALC!J17: 1163285 -> 926338.2
$ws.Cells.Item(17, 12).Value = 2779014.6  # ALC!L17: 3489855 -> 2779014.6
$ws.Cells.Item(17, 14).Value = -2779350.6  # ALC!N17: -3490191 -> -2779350.6

$ws.Cells.Item(129, 8).Value = 1073.6296  # ALC!H129: 1123.4186 -> 1073.6296
$ws.Cells.Item(129, 10).Value = 1073.6296  # ALC!J129: 1123.4186 -> 1073.6296
$ws.Cells.Item(129, 12).Value = 3220.8888  # ALC!L129: 3370.2558 -> 3220.8888
$ws.Cells.Item(129, 14).Value = -13220.8888  # ALC!N129: -13370.2558 -> -13220.8888

$ws.Cells.Item(132, 8).Value = 20669.04  # ALC!H132: 21073.824 -> 20669.04
$ws.Cells.Item(132, 9).Value = 23286.479  # ALC!I132: 23803.4 -> 23286.479
$ws.Cells.Item(132, 11).Value = 69859.43700000001  # ALC!K132: 71410.20000000001 -> 69859.43700000001
$ws.Cells.Item(132, 13).Value = -67329.43700000001  # ALC!M132: -68880.20000000001 -> -67329.43700000001

$ws.Cells.Item(135, 8).Value = 641.25  # ALC!H135: 406.975 -> 641.25
$ws.Cells.Item(135, 9).Value = 627.6316  # ALC!I135: 396.08334 -> 627.6316
$ws.Cells.Item(135, 10).Value = 900  # ALC!J135: 505 -> 900
$ws.Cells.Item(135, 11).Value = 5648.6844  # ALC!K135: 3564.75006 -> 5648.6844
$ws.Cells.Item(135, 12).Value = 8100  # ALC!L135: 4545 -> 8100
$ws.Cells.Item(135, 13).Value = -3113.6844  # ALC!M135: -1029.75006 -> -3113.6844
$ws.Cells.Item(135, 14).Value = -13170  # ALC!N135: -9615 -> -13170

$ws.Cells.Item(138, 8).Value = 2472.03  # ALC!H138: 2496.48 -> 2472.03
$ws.Cells.Item(138, 9).Value = 1211.7273  # ALC!I138: 1227.9846 -> 1211.7273
$ws.Cells.Item(138, 10).Value = 4918.5  # ALC!J138: 4852.2573 -> 4918.5
$ws.Cells.Item(138, 11).Value = 3635.1819  # ALC!K138: 3683.9538 -> 3635.1819
$ws.Cells.Item(138, 12).Value = 14755.5  # ALC!L138: 14556.7719 -> 14755.5
$ws.Cells.Item(138, 13).Value = 1504.8181  # ALC!M138: 1456.0462 -> 1504.8181
$ws.Cells.Item(138, 14).Value = -25035.5  # ALC!N138: -24836.7719 -> -25035.5

$ws.Cells.Item(141, 8).Value = 4990.8335  # ALC!H141: 5210.052 -> 4990.8335
$ws.Cells.Item(141, 9).Value = 1196.098  # ALC!I141: 1252.8043 -> 1196.098
$ws.Cells.Item(141, 10).Value = 26494.334  # ALC!J141: 20379.5 -> 26494.334
$ws.Cells.Item(141, 11).Value = 3588.294  # ALC!K141: 3758.4129 -> 3588.294
$ws.Cells.Item(141, 12).Value = 79483.00199999999  # ALC!L141: 61138.5 -> 79483.00199999999
$ws.Cells.Item(141, 13).Value = 1591.706  # ALC!M141: 1421.5871 -> 1591.706
$ws.Cells.Item(141, 14).Value = -89843.00199999999  # ALC!N141: -71498.5 -> -89843.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1395.525  # ARM!H61: 1115.24 -> 1395.525
$ws.Cells.Item(61, 9).Value = 741.45  # ARM!I61: 508.41934 -> 741.45
$ws.Cells.Item(61, 10).Value = 2049.6  # ARM!J61: 2105.3157 -> 2049.6
$ws.Cells.Item(61, 11).Value = 741.45  # ARM!K61: 508.41934 -> 741.45
$ws.Cells.Item(61, 12).Value = 2049.6  # ARM!L61: 2105.3157 -> 2049.6
$ws.Cells.Item(61, 13).Value = -529.45  # ARM!M61: -296.41934 -> -529.45
$ws.Cells.Item(61, 14).Value = -2473.6  # ARM!N61: -2529.3157 -> -2473.6

$ws.Cells.Item(74, 8).Value = 3440.0652  # ARM!H74: 3009.3396 -> 3440.0652
$ws.Cells.Item(74, 9).Value = 3624.5  # ARM!I74: 3132.2654 -> 3624.5
$ws.Cells.Item(74, 11).Value = 3624.5  # ARM!K74: 3132.2654 -> 3624.5
$ws.Cells.Item(74, 13).Value = -2750.5  # ARM!M74: -2258.2654 -> -2750.5

$ws.Cells.Item(77, 8).Value = 3440.0652  # ARM!H77: 3009.3396 -> 3440.0652
$ws.Cells.Item(77, 9).Value = 3624.5  # ARM!I77: 3132.2654 -> 3624.5
$ws.Cells.Item(77, 11).Value = 18122.5  # ARM!K77: 15661.327 -> 18122.5
$ws.Cells.Item(77, 13).Value = -13754.5  # ARM!M77: -11293.327 -> -13754.5

$ws.Cells.Item(136, 8).Value = 1395.525  # ARM!H136: 1115.24 -> 1395.525
$ws.Cells.Item(136, 9).Value = 741.45  # ARM!I136: 508.41934 -> 741.45
$ws.Cells.Item(136, 10).Value = 2049.6  # ARM!J136: 2105.3157 -> 2049.6
$ws.Cells.Item(136, 11).Value = 2224.35  # ARM!K136: 1525.25802 -> 2224.35
$ws.Cells.Item(136, 12).Value = 6148.799999999999  # ARM!L136: 6315.9471 -> 6148.799999999999
$ws.Cells.Item(136, 13).Value = 325.6499999999996  # ARM!M136: 1024.74198 -> 325.6499999999996
$ws.Cells.Item(136, 14).Value = -11248.8  # ARM!N136: -11415.9471 -> -11248.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2573.0557  # BSM!H86: 2678.4119 -> 2573.0557
$ws.Cells.Item(86, 9).Value = 2591  # BSM!I86: 3006.9092 -> 2591
$ws.Cells.Item(86, 10).Value = 2483.3333  # BSM!J86: 2076.1667 -> 2483.3333
$ws.Cells.Item(86, 11).Value = 2591  # BSM!K86: 3006.9092 -> 2591
$ws.Cells.Item(86, 12).Value = 2483.3333  # BSM!L86: 2076.1667 -> 2483.3333
$ws.Cells.Item(86, 13).Value = -1468  # BSM!M86: -1883.9092 -> -1468
$ws.Cells.Item(86, 14).Value = -4729.3333  # BSM!N86: -4322.1667 -> -4729.3333

$ws.Cells.Item(89, 8).Value = 2573.0557  # BSM!H89: 2678.4119 -> 2573.0557
$ws.Cells.Item(89, 9).Value = 2591  # BSM!I89: 3006.9092 -> 2591
$ws.Cells.Item(89, 10).Value = 2483.3333  # BSM!J89: 2076.1667 -> 2483.3333
$ws.Cells.Item(89, 11).Value = 12955  # BSM!K89: 15034.546 -> 12955
$ws.Cells.Item(89, 12).Value = 12416.6665  # BSM!L89: 10380.8335 -> 12416.6665
$ws.Cells.Item(89, 13).Value = -7339  # BSM!M89: -9418.546 -> -7339
$ws.Cells.Item(89, 14).Value = -23648.6665  # BSM!N89: -21612.8335 -> -23648.6665

$ws.Cells.Item(134, 8).Value = 1074.2727  # BSM!H134: 1396.3158 -> 1074.2727
$ws.Cells.Item(134, 9).Value = 829.04254  # BSM!I134: 1087.742 -> 829.04254
$ws.Cells.Item(134, 10).Value = 2515  # BSM!J134: 2762.8572 -> 2515
$ws.Cells.Item(134, 11).Value = 2487.12762  # BSM!K134: 3263.226 -> 2487.12762
$ws.Cells.Item(134, 12).Value = 7545  # BSM!L134: 8288.571599999999 -> 7545
$ws.Cells.Item(134, 13).Value = 47.87237999999979  # BSM!M134: -728.2259999999997 -> 47.87237999999979
$ws.Cells.Item(134, 14).Value = -12615  # BSM!N134: -13358.5716 -> -12615

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2256.8025  # CRP!H31: 2305.8718 -> 2256.8025
$ws.Cells.Item(31, 9).Value = 1354.5778  # CRP!I31: 1395.9767 -> 1354.5778
$ws.Cells.Item(31, 10).Value = 3384.5833  # CRP!J31: 3423.743 -> 3384.5833
$ws.Cells.Item(31, 11).Value = 1354.5778  # CRP!K31: 1395.9767 -> 1354.5778
$ws.Cells.Item(31, 12).Value = 3384.5833  # CRP!L31: 3423.743 -> 3384.5833
$ws.Cells.Item(31, 13).Value = -1059.5778  # CRP!M31: -1100.9767 -> -1059.5778
$ws.Cells.Item(31, 14).Value = -3974.5833  # CRP!N31: -4013.743 -> -3974.5833

$ws.Cells.Item(34, 8).Value = 2256.8025  # CRP!H34: 2305.8718 -> 2256.8025
$ws.Cells.Item(34, 9).Value = 1354.5778  # CRP!I34: 1395.9767 -> 1354.5778
$ws.Cells.Item(34, 10).Value = 3384.5833  # CRP!J34: 3423.743 -> 3384.5833
$ws.Cells.Item(34, 11).Value = 1354.5778  # CRP!K34: 1395.9767 -> 1354.5778
$ws.Cells.Item(34, 12).Value = 3384.5833  # CRP!L34: 3423.743 -> 3384.5833
$ws.Cells.Item(34, 13).Value = -1152.5778  # CRP!M34: -1193.9767 -> -1152.5778
$ws.Cells.Item(34, 14).Value = -3788.5833  # CRP!N34: -3827.743 -> -3788.5833

$ws.Cells.Item(37, 8).Value = 0  # CRP!H37: 8000 -> 0
$ws.Cells.Item(37, 9).Value = 0  # CRP!I37: 8000 -> 0
$ws.Cells.Item(37, 11).Value = 0  # CRP!K37: 8000 -> 0
$ws.Cells.Item(37, 13).Value = $null  # CRP!M37: -7893 -> (cleared)

$ws.Cells.Item(57, 8).Value = 0  # CRP!H57: 32500 -> 0
$ws.Cells.Item(57, 9).Value = 0  # CRP!I57: 15000 -> 0
$ws.Cells.Item(57, 10).Value = 0  # CRP!J57: 50000 -> 0
$ws.Cells.Item(57, 11).Value = 0  # CRP!K57: 15000 -> 0
$ws.Cells.Item(57, 12).Value = 0  # CRP!L57: 50000 -> 0
$ws.Cells.Item(57, 13).Value = $null  # CRP!M57: -14440 -> (cleared)
$ws.Cells.Item(57, 14).Value = $null  # CRP!N57: -51120 -> (cleared)

$ws.Cells.Item(58, 8).Value = 1393.8939  # CRP!H58: 1388.4054 -> 1393.8939
$ws.Cells.Item(58, 9).Value = 1084.6666  # CRP!I58: 1055.1666 -> 1084.6666
$ws.Cells.Item(58, 10).Value = 2445.2666  # CRP!J58: 2288.15 -> 2445.2666
$ws.Cells.Item(58, 11).Value = 1084.6666  # CRP!K58: 1055.1666 -> 1084.6666
$ws.Cells.Item(58, 12).Value = 2445.2666  # CRP!L58: 2288.15 -> 2445.2666
$ws.Cells.Item(58, 13).Value = -881.6666  # CRP!M58: -852.1666 -> -881.6666
$ws.Cells.Item(58, 14).Value = -2851.2666  # CRP!N58: -2694.15 -> -2851.2666

$ws.Cells.Item(99, 8).Value = 4014864.2  # CRP!H99: 11993.272 -> 4014864.2
$ws.Cells.Item(99, 9).Value = 10668217  # CRP!I99: 2532.4 -> 10668217
$ws.Cells.Item(99, 10).Value = 22852.8  # CRP!J99: 19877.334 -> 22852.8
$ws.Cells.Item(99, 11).Value = 10668217  # CRP!K99: 2532.4 -> 10668217
$ws.Cells.Item(99, 12).Value = 22852.8  # CRP!L99: 19877.334 -> 22852.8
$ws.Cells.Item(99, 13).Value = -10666719  # CRP!M99: -1034.4 -> -10666719
$ws.Cells.Item(99, 14).Value = -25848.8  # CRP!N99: -22873.334 -> -25848.8

$ws.Cells.Item(126, 8).Value = 4014864.2  # CRP!H126: 11993.272 -> 4014864.2
$ws.Cells.Item(126, 9).Value = 10668217  # CRP!I126: 2532.4 -> 10668217
$ws.Cells.Item(126, 10).Value = 22852.8  # CRP!J126: 19877.334 -> 22852.8
$ws.Cells.Item(126, 11).Value = 32004651  # CRP!K126: 7597.200000000001 -> 32004651
$ws.Cells.Item(126, 12).Value = 68558.39999999999  # CRP!L126: 59632.00199999999 -> 68558.39999999999
$ws.Cells.Item(126, 13).Value = -32002181  # CRP!M126: -5127.200000000001 -> -32002181
$ws.Cells.Item(126, 14).Value = -73498.39999999999  # CRP!N126: -64572.00199999999 -> -73498.39999999999

$ws.Cells.Item(136, 8).Value = 1393.8939  # CRP!H136: 1388.4054 -> 1393.8939
$ws.Cells.Item(136, 9).Value = 1084.6666  # CRP!I136: 1055.1666 -> 1084.6666
$ws.Cells.Item(136, 10).Value = 2445.2666  # CRP!J136: 2288.15 -> 2445.2666
$ws.Cells.Item(136, 11).Value = 3253.9998  # CRP!K136: 3165.4998 -> 3253.9998
$ws.Cells.Item(136, 12).Value = 7335.7998  # CRP!L136: 6864.450000000001 -> 7335.7998
$ws.Cells.Item(136, 13).Value = -703.9998000000001  # CRP!M136: -615.4998000000001 -> -703.9998000000001
$ws.Cells.Item(136, 14).Value = -12435.7998  # CRP!N136: -11964.45 -> -12435.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 10000  # CUL!H3: 11330 -> 10000
$ws.Cells.Item(3, 9).Value = 10000  # CUL!I3: 11330 -> 10000
$ws.Cells.Item(3, 11).Value = 30000  # CUL!K3: 33990 -> 30000
$ws.Cells.Item(3, 13).Value = -29888  # CUL!M3: -33878 -> -29888

$ws.Cells.Item(4, 8).Value = 1426.5834  # CUL!H4: 1320 -> 1426.5834
$ws.Cells.Item(4, 9).Value = 279.75  # CUL!I4: 280 -> 279.75
$ws.Cells.Item(4, 10).Value = 2000  # CUL!J4: 1840 -> 2000
$ws.Cells.Item(4, 11).Value = 839.25  # CUL!K4: 840 -> 839.25
$ws.Cells.Item(4, 12).Value = 6000  # CUL!L4: 5520 -> 6000
$ws.Cells.Item(4, 13).Value = -727.25  # CUL!M4: -728 -> -727.25
$ws.Cells.Item(4, 14).Value = -6224  # CUL!N4: -5744 -> -6224

$ws.Cells.Item(5, 8).Value = 651072  # CUL!H5: 676082.9 -> 651072
$ws.Cells.Item(5, 10).Value = 976081.1  # CUL!J5: 1033451.2 -> 976081.1
$ws.Cells.Item(5, 12).Value = 2928243.3  # CUL!L5: 3100353.6 -> 2928243.3
$ws.Cells.Item(5, 14).Value = -2928467.3  # CUL!N5: -3100577.6 -> -2928467.3

$ws.Cells.Item(22, 8).Value = 2070  # CUL!H22: 1286.9565 -> 2070
$ws.Cells.Item(22, 10).Value = 2916.6667  # CUL!J22: 2666.6667 -> 2916.6667
$ws.Cells.Item(22, 12).Value = 8750.000100000001  # CUL!L22: 8000.000100000001 -> 8750.000100000001
$ws.Cells.Item(22, 14).Value = -9088.000100000001  # CUL!N22: -8338.000100000001 -> -9088.000100000001

$ws.Cells.Item(27, 8).Value = 2070  # CUL!H27: 1286.9565 -> 2070
$ws.Cells.Item(27, 10).Value = 2916.6667  # CUL!J27: 2666.6667 -> 2916.6667
$ws.Cells.Item(27, 12).Value = 8750.000100000001  # CUL!L27: 8000.000100000001 -> 8750.000100000001
$ws.Cells.Item(27, 14).Value = -8954.000100000001  # CUL!N27: -8204.000100000001 -> -8954.000100000001

$ws.Cells.Item(34, 8).Value = 907.2857  # CUL!H34: 1070.2 -> 907.2857
$ws.Cells.Item(34, 10).Value = 1010.6  # CUL!J34: 1351 -> 1010.6
$ws.Cells.Item(34, 12).Value = 3031.8  # CUL!L34: 4053 -> 3031.8
$ws.Cells.Item(34, 14).Value = -3199.8  # CUL!N34: -4221 -> -3199.8

$ws.Cells.Item(113, 8).Value = 1036.5714  # CUL!H113: 1136.0555 -> 1036.5714
$ws.Cells.Item(113, 9).Value = 1125.875  # CUL!I113: 1169.9333 -> 1125.875
$ws.Cells.Item(113, 10).Value = 750.8  # CUL!J113: 966.6667 -> 750.8
$ws.Cells.Item(113, 11).Value = 3377.625  # CUL!K113: 3509.7999 -> 3377.625
$ws.Cells.Item(113, 12).Value = 2252.4  # CUL!L113: 2900.0001 -> 2252.4
$ws.Cells.Item(113, 13).Value = -1207.625  # CUL!M113: -1339.7999 -> -1207.625
$ws.Cells.Item(113, 14).Value = -6592.4  # CUL!N113: -7240.0001 -> -6592.4

$ws.Cells.Item(131, 8).Value = 3504.1025  # CUL!H131: 3975.7144 -> 3504.1025
$ws.Cells.Item(131, 9).Value = 375.29413  # CUL!I131: 610 -> 375.29413
$ws.Cells.Item(131, 10).Value = 5921.8184  # CUL!J131: 5322 -> 5921.8184
$ws.Cells.Item(131, 11).Value = 1125.88239  # CUL!K131: 1830 -> 1125.88239
$ws.Cells.Item(131, 12).Value = 17765.4552  # CUL!L131: 15966 -> 17765.4552
$ws.Cells.Item(131, 13).Value = 3914.11761  # CUL!M131: 3210 -> 3914.11761
$ws.Cells.Item(131, 14).Value = -27845.4552  # CUL!N131: -26046 -> -27845.4552

$ws.Cells.Item(135, 8).Value = 651072  # CUL!H135: 676082.9 -> 651072
$ws.Cells.Item(135, 10).Value = 976081.1  # CUL!J135: 1033451.2 -> 976081.1
$ws.Cells.Item(135, 12).Value = 8784729.9  # CUL!L135: 9301060.799999999 -> 8784729.9
$ws.Cells.Item(135, 14).Value = -8789799.9  # CUL!N135: -9306130.799999999 -> -8789799.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1021.3019  # GSM!H132: 1208.4773 -> 1021.3019
$ws.Cells.Item(132, 9).Value = 942.72  # GSM!I132: 1126.3414 -> 942.72
$ws.Cells.Item(132, 11).Value = 2828.16  # GSM!K132: 3379.0242 -> 2828.16
$ws.Cells.Item(132, 13).Value = -298.1599999999999  # GSM!M132: -849.0241999999998 -> -298.1599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 951.2632  # LTW!H16: 694.65515 -> 951.2632
$ws.Cells.Item(16, 9).Value = 576.5714  # LTW!I16: 424.0476 -> 576.5714
$ws.Cells.Item(16, 10).Value = 2000.4  # LTW!J16: 1405 -> 2000.4
$ws.Cells.Item(16, 11).Value = 576.5714  # LTW!K16: 424.0476 -> 576.5714
$ws.Cells.Item(16, 12).Value = 2000.4  # LTW!L16: 1405 -> 2000.4
$ws.Cells.Item(16, 13).Value = -406.5714  # LTW!M16: -254.0476 -> -406.5714
$ws.Cells.Item(16, 14).Value = -2340.4  # LTW!N16: -1745 -> -2340.4

$ws.Cells.Item(132, 8).Value = 1382.28  # LTW!H132: 4866.2856 -> 1382.28
$ws.Cells.Item(132, 9).Value = 1004.97437  # LTW!I132: 4433.797 -> 1004.97437
$ws.Cells.Item(132, 10).Value = 2720  # LTW!J132: 6250.25 -> 2720
$ws.Cells.Item(132, 11).Value = 3014.92311  # LTW!K132: 13301.391 -> 3014.92311
$ws.Cells.Item(132, 12).Value = 8160  # LTW!L132: 18750.75 -> 8160
$ws.Cells.Item(132, 13).Value = -484.9231100000002  # LTW!M132: -10771.391 -> -484.9231100000002
$ws.Cells.Item(132, 14).Value = -13220  # LTW!N132: -23810.75 -> -13220

$ws.Cells.Item(136, 8).Value = 11112984  # LTW!H136: 11496176 -> 11112984
$ws.Cells.Item(136, 9).Value = 2048.4814  # LTW!I136: 2105.9614 -> 2048.4814
$ws.Cells.Item(136, 10).Value = 111111400  # LTW!J136: 111111460 -> 111111400
$ws.Cells.Item(136, 11).Value = 6145.4442  # LTW!K136: 6317.8842 -> 6145.4442
$ws.Cells.Item(136, 12).Value = 333334200  # LTW!L136: 333334380 -> 333334200
$ws.Cells.Item(136, 13).Value = -3595.4442  # LTW!M136: -3767.8842 -> -3595.4442
$ws.Cells.Item(136, 14).Value = -333339300  # LTW!N136: -333339480 -> -333339300

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1452.0817  # WVR!H132: 1559.8572 -> 1452.0817
$ws.Cells.Item(132, 10).Value = 3001  # WVR!J132: 2600.4167 -> 3001
$ws.Cells.Item(132, 12).Value = 9003  # WVR!L132: 7801.250100000001 -> 9003
$ws.Cells.Item(132, 14).Value = -14063  # WVR!N132: -12861.2501 -> -14063

$ws.Cells.Item(136, 8).Value = 1369.1094  # WVR!H136: 1389.5238 -> 1369.1094
$ws.Cells.Item(136, 9).Value = 483  # WVR!I136: 490.9 -> 483
$ws.Cells.Item(136, 10).Value = 4845.385  # WVR!J136: 4845.769 -> 4845.385
$ws.Cells.Item(136, 11).Value = 1449  # WVR!K136: 1472.7 -> 1449
$ws.Cells.Item(136, 12).Value = 14536.155  # WVR!L136: 14537.307 -> 14536.155
$ws.Cells.Item(136, 13).Value = 1101  # WVR!M136: 1077.3 -> 1101
$ws.Cells.Item(136, 14).Value = -19636.155  # WVR!N136: -19637.307 -> -19636.155
